$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "it threw an exception."
$ws.Range("E4").WrapText = $true
$ws.Range("E4").VerticalAlignment = -4160

$ws.Range("E4").Select()
